# ---------------------------------------------------------------------------
# Add a new "localdb" command group to the #system sheet of the showcase
# workbook, per commit:
#   "NEW command type to create, maintain and manipulate a local-only
#    relational database ..."
#
# The #system sheet holds one column per command-group ("target"): a header
# cell with the group name, followed by the function signatures belonging to
# that group. Column A holds the alphabetically sorted list of all group
# names (the "target" defined name). Adding the "localdb" group means:
#   1. Inserting a new column right before the existing "macro" column (N),
#      which shifts every subsequent group's column one to the right.
#   2. Filling that new column with the "localdb" header and its six
#      function signatures (also alphabetically sorted).
#   3. Inserting "localdb" into the sorted list in column A (between "json"
#      and "macro"), shifting the later entries down by one row.
#   4. Re-pointing every defined name whose column moved, and adding the new
#      "localdb" defined name.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Insert a new, blank column before column N ("macro"). This shifts
#    columns N..AD (macro, mail, number, pdf, rdbms, redis, sms, sound, ssh,
#    step, web, webalert, webcookie, ws, ws.async, xml) one column to the
#    right, into O..AE.
# ---------------------------------------------------------------------------
$ws.Columns("N").Insert()

# ---------------------------------------------------------------------------
# 2) Populate the new column N with the "localdb" group: header + the six
#    function signatures, alphabetically ordered.
# ---------------------------------------------------------------------------
$ws.Range("N1").Value = "localdb"
$ws.Range("N2").Value = "cloneTable(var,source,target)"
$ws.Range("N3").Value = "dropTables(var,tables)"
$ws.Range("N4").Value = "exportCSV(sql,output)"
$ws.Range("N5").Value = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value = "purge(var)"
$ws.Range("N7").Value = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------------
# 3) Insert "localdb" into the alphabetically sorted group list in column A.
#    "localdb" sorts right after "json" (row 13) and before "macro"
#    (row 14), so shift A14:A29 down to A15:A30, then drop "localdb" into
#    the newly vacated A14.
# ---------------------------------------------------------------------------
for ($r = 29; $r -ge 14; $r--) {
    $ws.Range("A$($r + 1)").Value = $ws.Range("A$r").Value()
}
$ws.Range("A14").Value = "localdb"

# ---------------------------------------------------------------------------
# 4) Fix up defined names: every group whose column shifted right needs its
#    RefersTo updated, "target" now spans one more row, and "localdb" is a
#    brand new named range.
# ---------------------------------------------------------------------------
$wb.Names.Item("mail").RefersTo      = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo    = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo       = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo     = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo     = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo       = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo     = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo       = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo      = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo       = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo  = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo     = "='#system'!`$O`$2:`$O`$4"

$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
